$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.718.82'
$ws.Range("E2").Value = '  +0.27%  '

$ws.Range("D3").Value = '2.622.14'
$ws.Range("E3").Value = '  -1.01%  '

$ws.Range("E4").Value = '  +0.06%  '

$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '595.68'
$cell.Style = $origStyle
$ws.Range("E5").Value = '  -1.05%  '

$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '149.91'
$cell.Style = $origStyle
$ws.Range("E6").Value = '  +2.01%  '

$ws.Range("E7").Value = '  +0.03%  '

$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.588'
$cell.Style = $origStyle
$ws.Range("E8").Value = '  -0.21%  '

$ws.Range("E9").Value = '  +0.36%  '

$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.68'
$cell.Style = $origStyle
$ws.Range("E10").Value = '  +1.71%  '

$ws.Range("E11").Value = '  +3.28%  '

$ws.Range("E12").Value = '  -1.16%  '

$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '27.65'
$cell.Style = $origStyle
$ws.Range("E13").Value = '  +0.58%  '

$ws.Range("D14").Value = '3.095.61'
$ws.Range("E14").Value = '  -0.98%  '

$ws.Range("D15").Value = '63.583.42'
$ws.Range("E15").Value = '  +0.25%  '

$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0000149'
$cell.Style = $origStyle
$ws.Range("E16").Value = '  +2.31%  '

$ws.Range("D17").Value = '2.626.80'
$ws.Range("E17").Value = '  -1.36%  '

$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '12.27'
$cell.Style = $origStyle
$ws.Range("E18").Value = '  +7.07%  '

$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '4.63'
$cell.Style = $origStyle
$ws.Range("E19").Value = '  +1.73%  '

$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '348.22'
$cell.Style = $origStyle
$ws.Range("E20").Value = '  +1.90%  '

$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '6.85'
$cell.Style = $origStyle
$ws.Range("E21").Value = '  -1.34%  '

$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.998'
$cell.Style = $origStyle
$ws.Range("E22").Value = '  -0.25%  '

$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.70'
$cell.Style = $origStyle
$ws.Range("E23").Value = '  +1.95%  '

$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '66.28'
$cell.Style = $origStyle
$ws.Range("E24").Value = '  -0.63%  '

$ws.Range("E25").Value = '  +11.42%  '

$cell = $ws.Range("D26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '9.18'
$cell.Style = $origStyle
$ws.Range("E26").Value = '  +1.09%  '

$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.66'
$cell.Style = $origStyle
$ws.Range("E27").Value = '  -1.47%  '

$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '561.15'
$cell.Style = $origStyle
$ws.Range("E28").Value = '  +0.18%  '

$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '8.22'
$cell.Style = $origStyle
$ws.Range("E29").Value = '  +3.60%  '

$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.164'
$cell.Style = $origStyle
$ws.Range("E30").Value = '  +0.27%  '

$ws.Range("E31").Value = '  -0.05%  '

$ws.Range("E32").Value = '  +0.34%  '

$ws.Range("D33").Value = '0.0₃0840'
$ws.Range("E33").Value = '  +2.78%  '

$ws.Range("E34").Value = '  -0.28%  '

$cell = $ws.Range("D35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.20'
$cell.Style = $origStyle
$ws.Range("E35").Value = '  +0.71%  '

$cell = $ws.Range("D36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '168.76'
$cell.Style = $origStyle
$ws.Range("E36").Value = '  +0.72%  '

$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.408'
$cell.Style = $origStyle
$ws.Range("E37").Value = '  +0.23%  '

$ws.Range("E38").Value = '  -0.05%  '

$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.93'
$cell.Style = $origStyle
$ws.Range("E39").Value = '  +0.47%  '

$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '19.30'
$cell.Style = $origStyle
$ws.Range("E40").Value = '  +0.99%  '

$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = $origStyle
$ws.Range("E41").Value = '  -0.05%  '

$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '169.36'
$cell.Style = $origStyle
$ws.Range("E42").Value = '  +0.53%  '

$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '39.90'
$cell.Style = $origStyle
$ws.Range("E43").Value = '  -0.07%  '

$ws.Range("E44").Value = '  +3.71%  '

$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0594'
$cell.Style = $origStyle
$ws.Range("E45").Value = '  +4.02%  '

$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '21.28'
$cell.Style = $origStyle
$ws.Range("E46").Value = '  -3.93%  '

$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.629'
$cell.Style = $origStyle
$ws.Range("E47").Value = '  -0.10%  '

$ws.Range("E48").Value = '  +0.42%  '

$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.97'
$cell.Style = $origStyle
$ws.Range("E49").Value = '  +5.37%  '

$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0968'
$cell.Style = $origStyle
$ws.Range("E50").Value = '  +0.84%  '

$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '19.12'
$cell.Style = $origStyle
$ws.Range("E51").Value = '  +1.77%  '
